$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (46060 -> 46061, i.e. 2026-02-07 -> 2026-02-08) for every data row (2..146).
$ws.Range("C2:C146").Value = 46061
